$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Combine the data from rows 2-7 into a single cell in A2, matching the
# Python tuple/list repr format, then remove the now-empty rows 3-7.
$ws.Range("A2").Value = "('Helm of Kaldra', ['{3}', 'Legendary Artifact — Equipment', 'Equipped creature has first strike, trample, and haste.', '{1}: If you control Equipment named Helm of Kaldra, Sword of Kaldra, and Shield of Kaldra, create Kaldra, a legendary 4/4 colorless Avatar creature token. Attach those Equipment to it.', 'Equip {2}'])"

$ws.Range("A3:A7").EntireRow.Delete()
